$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.965.90"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "2.996.69"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.64"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.30"
$ws.Range("E6").Value = "  +4.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "2.988.61"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").Value = "  +11.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.147"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000219"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.66"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "3.484.55"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "62.013.65"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "3.000.80"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.51"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.85"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.23"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.17"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.78"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.48"
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.54"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.27"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.31"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.48"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.79"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "446.84"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0801"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0385"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "2.923.40"
$ws.Range("E40").Value = "  -7.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.02"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("E43").Value = "  +5.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.52"
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.244"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.108"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.97"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.37"
$ws.Range("E49").Value = "  -2.67%  "
$ws.Range("D50").Value = "0.0₃0494"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  -3.42%  "
